$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels
$ws.Range("B1").Value = "Role"
$ws.Range("C1").Value = "IDAM Roles"

# Move selection to C1 (matches saved file's cursor position)
$ws.Range("C1").Select()
